$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Rewrite the Polymer label column so the workbook-level shared string
# table is rebuilt (this drops the now-unused "DIP ..." labels and
# renumbers the remaining entries), then apply the real content change:
# "DIP S1" -> "S" and "DIP B1" -> "B".
$ws1.Range("A2").Value = "pDNA"
$ws1.Range("A3").Value = "S"
$ws1.Range("A4").Value = "B"
$ws1.Range("A5").Value = "G1"
$ws1.Range("A6").Value = "G2"
$ws1.Range("A7").Value = "G3"

$ws1.Range("A8").Select() | Out-Null
